$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"44.8529195"
$ws.Range("H2").Value = [double]"89.705839"
$ws.Range("I2").Value = [double]"0.09328277884630309"
$ws.Range("J2").Value = [double]"0.07358468181289031"
$ws.Range("M2").Value = [double]"10.836393"
$ws.Range("N2").Value = [double]"21.672786"
$ws.Range("O2").Value = [double]"0.01846862257356514"
$ws.Range("P2").Value = [double]"0.01262296693390161"
$ws.Range("Q2").Value = [double]"486.0438628993635"
$ws.Range("R2").Value = [double]"1944.175451597454"
$ws.Range("S2").Value = [double]"0.001722804435125718"
$ws.Range("T2").Value = [double]"0.0009288570053657855"
# Row 3
$ws.Range("G3").Value = [double]"44.8529195"
$ws.Range("H3").Value = [double]"89.705839"
$ws.Range("I3").Value = [double]"0.09328277884630309"
$ws.Range("J3").Value = [double]"0.07358468181289031"
$ws.Range("O3").Value = [double]"0.2380839126543345"
$ws.Range("P3").Value = [double]"0.2440890227431923"
$ws.Range("Q3").Value = [double]"6265.720366517257"
$ws.Range("R3").Value = [double]"37594.32219910355"
$ws.Range("S3").Value = [double]"0.02220912897099683"
$ws.Range("T3").Value = [double]"0.01796121307257715"
# Row 4
$ws.Range("G4").Value = [double]"44.8529195"
$ws.Range("H4").Value = [double]"89.705839"
$ws.Range("I4").Value = [double]"0.09328277884630309"
$ws.Range("J4").Value = [double]"0.07358468181289031"
$ws.Range("M4").Value = [double]"150.12088"
$ws.Range("N4").Value = [double]"450.36264"
$ws.Range("O4").Value = [double]"0.2558532043948076"
$ws.Range("P4").Value = [double]"0.2623065033256284"
$ws.Range("Q4").Value = [double]"6733.35974590916"
$ws.Range("R4").Value = [double]"40400.15847545496"
$ws.Range("S4").Value = [double]"0.02386669788267882"
$ws.Range("T4").Value = [double]"0.01930174058466822"
# Row 5
$ws.Range("G5").Value = [double]"44.8529195"
$ws.Range("H5").Value = [double]"89.705839"
$ws.Range("I5").Value = [double]"0.09328277884630309"
$ws.Range("J5").Value = [double]"0.07358468181289031"
$ws.Range("M5").Value = [double]"32.469223"
$ws.Range("N5").Value = [double]"64.938446"
$ws.Range("O5").Value = [double]"0.05533777012737728"
$ws.Range("P5").Value = [double]"0.03782235733776705"
$ws.Range("Q5").Value = [double]"1456.339445446548"
$ws.Range("R5").Value = [double]"5825.357781786193"
$ws.Range("S5").Value = [double]"0.005162060972639692"
$ws.Range("T5").Value = [double]"0.002783146130113026"
# Row 6
$ws.Range("G6").Value = [double]"44.8529195"
$ws.Range("H6").Value = [double]"89.705839"
$ws.Range("I6").Value = [double]"0.09328277884630309"
$ws.Range("J6").Value = [double]"0.07358468181289031"
$ws.Range("M6").Value = [double]"52.09024033333333"
$ws.Range("N6").Value = [double]"156.270721"
$ws.Range("O6").Value = [double]"0.08877815602319267"
$ws.Range("P6").Value = [double]"0.09101737745760805"
$ws.Range("Q6").Value = [double]"2336.399356406653"
$ws.Range("R6").Value = [double]"14018.39613843992"
$ws.Range("S6").Value = [double]"0.008281473094694073"
$ws.Range("T6").Value = [double]"0.006697484759661823"
# Row 7
$ws.Range("G7").Value = [double]"44.8529195"
$ws.Range("H7").Value = [double]"89.705839"
$ws.Range("I7").Value = [double]"0.09328277884630309"
$ws.Range("J7").Value = [double]"0.07358468181289031"
$ws.Range("M7").Value = [double]"201.5345866666667"
$ws.Range("N7").Value = [double]"604.60376"
$ws.Range("O7").Value = [double]"0.3434783342267227"
$ws.Range("P7").Value = [double]"0.3521417722019025"
$ws.Range("Q7").Value = [double]"9039.414592225772"
$ws.Range("R7").Value = [double]"54236.48755335464"
$ws.Range("S7").Value = [double]"0.03204061349016796"
$ws.Range("T7").Value = [double]"0.0259122402605043"
# Row 8
$ws.Range("I8").Value = [double]"0.5330899821806619"
$ws.Range("J8").Value = [double]"0.630779719497357"
$ws.Range("M8").Value = [double]"10.836393"
$ws.Range("N8").Value = [double]"21.672786"
$ws.Range("O8").Value = [double]"0.01846862257356514"
$ws.Range("P8").Value = [double]"0.01262296693390161"
$ws.Range("Q8").Value = [double]"2777.630741885971"
$ws.Range("R8").Value = [double]"16665.78445131582"
$ws.Range("S8").Value = [double]"0.00984543767864321"
$ws.Range("T8").Value = [double]"0.007962311541790869"
# Row 9
$ws.Range("I9").Value = [double]"0.5330899821806619"
$ws.Range("J9").Value = [double]"0.630779719497357"
$ws.Range("O9").Value = [double]"0.2380839126543345"
$ws.Range("P9").Value = [double]"0.2440890227431923"
$ws.Range("S9").Value = [double]"0.1269201487544015"
$ws.Range("T9").Value = [double]"0.1539664052983349"
# Row 10
$ws.Range("I10").Value = [double]"0.5330899821806619"
$ws.Range("J10").Value = [double]"0.630779719497357"
$ws.Range("M10").Value = [double]"150.12088"
$ws.Range("N10").Value = [double]"450.36264"
$ws.Range("O10").Value = [double]"0.2558532043948076"
$ws.Range("P10").Value = [double]"0.2623065033256284"
$ws.Range("Q10").Value = [double]"38479.6279801752"
$ws.Range("R10").Value = [double]"346316.6518215768"
$ws.Range("S10").Value = [double]"0.1363927801716933"
$ws.Range("T10").Value = [double]"0.1654576225900725"
# Row 11
$ws.Range("I11").Value = [double]"0.5330899821806619"
$ws.Range("J11").Value = [double]"0.630779719497357"
$ws.Range("M11").Value = [double]"32.469223"
$ws.Range("N11").Value = [double]"64.938446"
$ws.Range("O11").Value = [double]"0.05533777012737728"
$ws.Range("P11").Value = [double]"0.03782235733776705"
$ws.Range("Q11").Value = [double]"8322.65053232667"
$ws.Range("R11").Value = [double]"49935.90319396002"
$ws.Range("S11").Value = [double]"0.02950001089112112"
$ws.Range("T11").Value = [double]"0.0238575759522455"
# Row 12
$ws.Range("I12").Value = [double]"0.5330899821806619"
$ws.Range("J12").Value = [double]"0.630779719497357"
$ws.Range("M12").Value = [double]"52.09024033333333"
$ws.Range("N12").Value = [double]"156.270721"
$ws.Range("O12").Value = [double]"0.08877815602319267"
$ws.Range("P12").Value = [double]"0.09101737745760805"
$ws.Range("Q12").Value = [double]"13351.99386937103"
$ws.Range("R12").Value = [double]"120167.9448243393"
$ws.Range("S12").Value = [double]"0.04732674561243581"
$ws.Range("T12").Value = [double]"0.05741191582209507"
# Row 13
$ws.Range("I13").Value = [double]"0.5330899821806619"
$ws.Range("J13").Value = [double]"0.630779719497357"
$ws.Range("M13").Value = [double]"201.5345866666667"
$ws.Range("N13").Value = [double]"604.60376"
$ws.Range("O13").Value = [double]"0.3434783342267227"
$ws.Range("P13").Value = [double]"0.3521417722019025"
$ws.Range("Q13").Value = [double]"51658.2098377768"
$ws.Range("R13").Value = [double]"464923.8885399912"
$ws.Range("S13").Value = [double]"0.1831048590723671"
$ws.Range("T13").Value = [double]"0.2221238882928183"
# Row 14
$ws.Range("G14").Value = [double]"0.08220233333333334"
$ws.Range("H14").Value = [double]"0.246607"
$ws.Range("I14").Value = [double]"0.0001709601552466038"
$ws.Range("J14").Value = [double]"0.000202288923777096"
$ws.Range("M14").Value = [double]"10.836393"
$ws.Range("N14").Value = [double]"21.672786"
$ws.Range("O14").Value = [double]"0.01846862257356514"
$ws.Range("P14").Value = [double]"0.01262296693390161"
$ws.Range("Q14").Value = [double]"0.8907767895170001"
$ws.Range("R14").Value = [double]"5.344660737102001"
$ws.Range("S14").Value = [double]"3.157398582367627E-06"
$ws.Range("T14").Value = [double]"2.553486395932825E-06"
# Row 15
$ws.Range("G15").Value = [double]"0.08220233333333334"
$ws.Range("H15").Value = [double]"0.246607"
$ws.Range("I15").Value = [double]"0.0001709601552466038"
$ws.Range("J15").Value = [double]"0.000202288923777096"
$ws.Range("O15").Value = [double]"0.2380839126543345"
$ws.Range("P15").Value = [double]"0.2440890227431923"
$ws.Range("Q15").Value = [double]"11.48323988457222"
$ws.Range("R15").Value = [double]"103.34915896115"
$ws.Range("S15").Value = [double]"4.070286266910389E-05"
$ws.Range("T15").Value = [double]"4.937650571652348E-05"
# Row 16
$ws.Range("G16").Value = [double]"0.08220233333333334"
$ws.Range("H16").Value = [double]"0.246607"
$ws.Range("I16").Value = [double]"0.0001709601552466038"
$ws.Range("J16").Value = [double]"0.000202288923777096"
$ws.Range("M16").Value = [double]"150.12088"
$ws.Range("N16").Value = [double]"450.36264"
$ws.Range("O16").Value = [double]"0.2558532043948076"
$ws.Range("P16").Value = [double]"0.2623065033256284"
$ws.Range("Q16").Value = [double]"12.34028661805333"
$ws.Range("R16").Value = [double]"111.06257956248"
$ws.Range("S16").Value = [double]"4.374070354367737E-05"
$ws.Range("T16").Value = [double]"5.306170025747462E-05"
# Row 17
$ws.Range("G17").Value = [double]"0.08220233333333334"
$ws.Range("H17").Value = [double]"0.246607"
$ws.Range("I17").Value = [double]"0.0001709601552466038"
$ws.Range("J17").Value = [double]"0.000202288923777096"
$ws.Range("M17").Value = [double]"32.469223"
$ws.Range("N17").Value = [double]"64.938446"
$ws.Range("O17").Value = [double]"0.05533777012737728"
$ws.Range("P17").Value = [double]"0.03782235733776705"
$ws.Range("Q17").Value = [double]"2.669045892120333"
$ws.Range("R17").Value = [double]"16.014275352722"
$ws.Range("S17").Value = [double]"9.460553771977294E-06"
$ws.Range("T17").Value = [double]"7.651043960569644E-06"
# Row 18
$ws.Range("G18").Value = [double]"0.08220233333333334"
$ws.Range("H18").Value = [double]"0.246607"
$ws.Range("I18").Value = [double]"0.0001709601552466038"
$ws.Range("J18").Value = [double]"0.000202288923777096"
$ws.Range("M18").Value = [double]"52.09024033333333"
$ws.Range("N18").Value = [double]"156.270721"
$ws.Range("O18").Value = [double]"0.08877815602319267"
$ws.Range("P18").Value = [double]"0.09101737745760805"
$ws.Range("Q18").Value = [double]"4.281939299294111"
$ws.Range("R18").Value = [double]"38.537453693647"
$ws.Range("S18").Value = [double]"1.517752733623223E-05"
$ws.Range("T18").Value = [double]"1.841180733091325E-05"
# Row 19
$ws.Range("G19").Value = [double]"0.08220233333333334"
$ws.Range("H19").Value = [double]"0.246607"
$ws.Range("I19").Value = [double]"0.0001709601552466038"
$ws.Range("J19").Value = [double]"0.000202288923777096"
$ws.Range("M19").Value = [double]"201.5345866666667"
$ws.Range("N19").Value = [double]"604.60376"
$ws.Range("O19").Value = [double]"0.3434783342267227"
$ws.Range("P19").Value = [double]"0.3521417722019025"
$ws.Range("Q19").Value = [double]"16.56661327136889"
$ws.Range("R19").Value = [double]"149.09951944232"
$ws.Range("S19").Value = [double]"5.872110934324538E-05"
$ws.Range("T19").Value = [double]"7.123438011568215E-05"
# Row 20
$ws.Range("G20").Value = [double]"178.5463335"
$ws.Range("H20").Value = [double]"357.092667"
$ws.Range("I20").Value = [double]"0.3713314167141066"
$ws.Range("J20").Value = [double]"0.2929190627035035"
$ws.Range("M20").Value = [double]"10.836393"
$ws.Range("N20").Value = [double]"21.672786"
$ws.Range("O20").Value = [double]"0.01846862257356514"
$ws.Range("P20").Value = [double]"0.01262296693390161"
$ws.Range("Q20").Value = [double]"1934.798238515066"
$ws.Range("R20").Value = [double]"7739.192954060263"
$ws.Range("S20").Value = [double]"0.006857979785000072"
$ws.Range("T20").Value = [double]"0.003697507642815777"
# Row 21
$ws.Range("G21").Value = [double]"178.5463335"
$ws.Range("H21").Value = [double]"357.092667"
$ws.Range("I21").Value = [double]"0.3713314167141066"
$ws.Range("J21").Value = [double]"0.2929190627035035"
$ws.Range("O21").Value = [double]"0.2380839126543345"
$ws.Range("P21").Value = [double]"0.2440890227431923"
$ws.Range("Q21").Value = [double]"24941.99732478802"
$ws.Range("R21").Value = [double]"149651.9839487281"
$ws.Range("S21").Value = [double]"0.08840803658277166"
$ws.Range("T21").Value = [double]"0.07149832775815006"
# Row 22
$ws.Range("G22").Value = [double]"178.5463335"
$ws.Range("H22").Value = [double]"357.092667"
$ws.Range("I22").Value = [double]"0.3713314167141066"
$ws.Range("J22").Value = [double]"0.2929190627035035"
$ws.Range("M22").Value = [double]"150.12088"
$ws.Range("N22").Value = [double]"450.36264"
$ws.Range("O22").Value = [double]"0.2558532043948076"
$ws.Range("P22").Value = [double]"0.2623065033256284"
$ws.Range("Q22").Value = [double]"26803.53270579348"
$ws.Range("R22").Value = [double]"160821.1962347609"
$ws.Range("S22").Value = [double]"0.09500633285876781"
$ws.Range("T22").Value = [double]"0.07683457509517652"
# Row 23
$ws.Range("G23").Value = [double]"178.5463335"
$ws.Range("H23").Value = [double]"357.092667"
$ws.Range("I23").Value = [double]"0.3713314167141066"
$ws.Range("J23").Value = [double]"0.2929190627035035"
$ws.Range("M23").Value = [double]"32.469223"
$ws.Range("N23").Value = [double]"64.938446"
$ws.Range("O23").Value = [double]"0.05533777012737728"
$ws.Range("P23").Value = [double]"0.03782235733776705"
$ws.Range("Q23").Value = [double]"5797.260718243871"
$ws.Range("R23").Value = [double]"23189.04287297548"
$ws.Range("S23").Value = [double]"0.02054865257919857"
$ws.Range("T23").Value = [double]"0.0110788894606157"
# Row 24
$ws.Range("G24").Value = [double]"178.5463335"
$ws.Range("H24").Value = [double]"357.092667"
$ws.Range("I24").Value = [double]"0.3713314167141066"
$ws.Range("J24").Value = [double]"0.2929190627035035"
$ws.Range("M24").Value = [double]"52.09024033333333"
$ws.Range("N24").Value = [double]"156.270721"
$ws.Range("O24").Value = [double]"0.08877815602319267"
$ws.Range("P24").Value = [double]"0.09101737745760805"
$ws.Range("Q24").Value = [double]"9300.521422650483"
$ws.Range("R24").Value = [double]"55803.1285359029"
$ws.Range("S24").Value = [double]"0.03296611844935813"
$ws.Range("T24").Value = [double]"0.02666072489461354"
# Row 25
$ws.Range("G25").Value = [double]"178.5463335"
$ws.Range("H25").Value = [double]"357.092667"
$ws.Range("I25").Value = [double]"0.3713314167141066"
$ws.Range("J25").Value = [double]"0.2929190627035035"
$ws.Range("M25").Value = [double]"201.5345866666667"
$ws.Range("N25").Value = [double]"604.60376"
$ws.Range("O25").Value = [double]"0.3434783342267227"
$ws.Range("P25").Value = [double]"0.3521417722019025"
$ws.Range("Q25").Value = [double]"35983.26152277132"
$ws.Range("R25").Value = [double]"215899.5691366279"
$ws.Range("S25").Value = [double]"0.1275442964590104"
$ws.Range("T25").Value = [double]"0.1031490378521319"
# Row 26
$ws.Range("E26").Value = [double]"2"
$ws.Range("F26").Value = [double]"0.6666666666666666"
$ws.Range("G26").Value = [double]"0.7601923333333334"
$ws.Range("H26").Value = [double]"2.280577"
$ws.Range("I26").Value = [double]"0.00158100864116523"
$ws.Range("J26").Value = [double]"0.001870731434715147"
$ws.Range("M26").Value = [double]"10.836393"
$ws.Range("N26").Value = [double]"21.672786"
$ws.Range("O26").Value = [double]"0.01846862257356514"
$ws.Range("P26").Value = [double]"0.01262296693390161"
$ws.Range("Q26").Value = [double]"8.237742879587001"
$ws.Range("R26").Value = [double]"49.426457277522"
$ws.Range("S26").Value = [double]"2.919905187922572E-05"
$ws.Range("T26").Value = [double]"2.361418104261961E-05"
# Row 27
$ws.Range("E27").Value = [double]"2"
$ws.Range("F27").Value = [double]"0.6666666666666666"
$ws.Range("G27").Value = [double]"0.7601923333333334"
$ws.Range("H27").Value = [double]"2.280577"
$ws.Range("I27").Value = [double]"0.00158100864116523"
$ws.Range("J27").Value = [double]"0.001870731434715147"
$ws.Range("O27").Value = [double]"0.2380839126543345"
$ws.Range("P27").Value = [double]"0.2440890227431923"
$ws.Range("Q27").Value = [double]"106.1949286364055"
$ws.Range("R27").Value = [double]"955.7543577276499"
$ws.Range("S27").Value = [double]"0.0003764127232289309"
$ws.Range("T27").Value = [double]"0.0004566250077145902"
# Row 28
$ws.Range("E28").Value = [double]"2"
$ws.Range("F28").Value = [double]"0.6666666666666666"
$ws.Range("G28").Value = [double]"0.7601923333333334"
$ws.Range("H28").Value = [double]"2.280577"
$ws.Range("I28").Value = [double]"0.00158100864116523"
$ws.Range("J28").Value = [double]"0.001870731434715147"
$ws.Range("M28").Value = [double]"150.12088"
$ws.Range("N28").Value = [double]"450.36264"
$ws.Range("O28").Value = [double]"0.2558532043948076"
$ws.Range("P28").Value = [double]"0.2623065033256284"
$ws.Range("Q28").Value = [double]"114.1207420492533"
$ws.Range("R28").Value = [double]"1027.08667844328"
$ws.Range("S28").Value = [double]"0.0004045061270180047"
$ws.Range("T28").Value = [double]"0.0004907050213014663"
# Row 29
$ws.Range("E29").Value = [double]"2"
$ws.Range("F29").Value = [double]"0.6666666666666666"
$ws.Range("G29").Value = [double]"0.7601923333333334"
$ws.Range("H29").Value = [double]"2.280577"
$ws.Range("I29").Value = [double]"0.00158100864116523"
$ws.Range("J29").Value = [double]"0.001870731434715147"
$ws.Range("M29").Value = [double]"32.469223"
$ws.Range("N29").Value = [double]"64.938446"
$ws.Range("O29").Value = [double]"0.05533777012737728"
$ws.Range("P29").Value = [double]"0.03782235733776705"
$ws.Range("Q29").Value = [double]"24.68285439389033"
$ws.Range("R29").Value = [double]"148.097126363342"
$ws.Range("S29").Value = [double]"8.748949275419863E-05"
$ws.Range("T29").Value = [double]"7.07554728067899E-05"
# Row 30
$ws.Range("E30").Value = [double]"2"
$ws.Range("F30").Value = [double]"0.6666666666666666"
$ws.Range("G30").Value = [double]"0.7601923333333334"
$ws.Range("H30").Value = [double]"2.280577"
$ws.Range("I30").Value = [double]"0.00158100864116523"
$ws.Range("J30").Value = [double]"0.001870731434715147"
$ws.Range("M30").Value = [double]"52.09024033333333"
$ws.Range("N30").Value = [double]"156.270721"
$ws.Range("O30").Value = [double]"0.08877815602319267"
$ws.Range("P30").Value = [double]"0.09101737745760805"
$ws.Range("Q30").Value = [double]"39.59860134289077"
$ws.Range("R30").Value = [double]"356.387412086017"
$ws.Range("S30").Value = [double]"0.0001403590318193827"
$ws.Range("T30").Value = [double]"0.0001702690691152811"
# Row 31
$ws.Range("E31").Value = [double]"2"
$ws.Range("F31").Value = [double]"0.6666666666666666"
$ws.Range("G31").Value = [double]"0.7601923333333334"
$ws.Range("H31").Value = [double]"2.280577"
$ws.Range("I31").Value = [double]"0.00158100864116523"
$ws.Range("J31").Value = [double]"0.001870731434715147"
$ws.Range("M31").Value = [double]"201.5345866666667"
$ws.Range("N31").Value = [double]"604.60376"
$ws.Range("O31").Value = [double]"0.3434783342267227"
$ws.Range("P31").Value = [double]"0.3521417722019025"
$ws.Range("Q31").Value = [double]"153.2050476855022"
$ws.Range("R31").Value = [double]"1378.84542916952"
$ws.Range("S31").Value = [double]"0.0005430422144654878"
$ws.Range("T31").Value = [double]"0.0006587626827343994"
# Row 32
$ws.Range("G32").Value = [double]"0.2614996666666667"
$ws.Range("H32").Value = [double]"0.7844990000000001"
$ws.Range("I32").Value = [double]"0.0005438534625164957"
$ws.Range("J32").Value = [double]"0.0006435156277567465"
$ws.Range("M32").Value = [double]"10.836393"
$ws.Range("N32").Value = [double]"21.672786"
$ws.Range("O32").Value = [double]"0.01846862257356514"
$ws.Range("P32").Value = [double]"0.01262296693390161"
$ws.Range("Q32").Value = [double]"2.833713157369"
$ws.Range("R32").Value = [double]"17.002278944214"
$ws.Range("S32").Value = [double]"1.004422433454371E-05"
$ws.Range("T32").Value = [double]"8.123076490622347E-06"
# Row 33
$ws.Range("G33").Value = [double]"0.2614996666666667"
$ws.Range("H33").Value = [double]"0.7844990000000001"
$ws.Range("I33").Value = [double]"0.0005438534625164957"
$ws.Range("J33").Value = [double]"0.0006435156277567465"
$ws.Range("O33").Value = [double]"0.2380839126543345"
$ws.Range("P33").Value = [double]"0.2440890227431923"
$ws.Range("Q33").Value = [double]"36.53014799339444"
$ws.Range("R33").Value = [double]"328.77133194055"
$ws.Range("S33").Value = [double]"0.0001294827602665347"
$ws.Range("T33").Value = [double]"0.0001570751006991162"
# Row 34
$ws.Range("G34").Value = [double]"0.2614996666666667"
$ws.Range("H34").Value = [double]"0.7844990000000001"
$ws.Range("I34").Value = [double]"0.0005438534625164957"
$ws.Range("J34").Value = [double]"0.0006435156277567465"
$ws.Range("M34").Value = [double]"150.12088"
$ws.Range("N34").Value = [double]"450.36264"
$ws.Range("O34").Value = [double]"0.2558532043948076"
$ws.Range("P34").Value = [double]"0.2623065033256284"
$ws.Range("Q34").Value = [double]"39.25656007970667"
$ws.Range("R34").Value = [double]"353.30904071736"
$ws.Range("S34").Value = [double]"0.0001391466511060568"
$ws.Range("T34").Value = [double]"0.0001687983341522689"
# Row 35
$ws.Range("G35").Value = [double]"0.2614996666666667"
$ws.Range("H35").Value = [double]"0.7844990000000001"
$ws.Range("I35").Value = [double]"0.0005438534625164957"
$ws.Range("J35").Value = [double]"0.0006435156277567465"
$ws.Range("M35").Value = [double]"32.469223"
$ws.Range("N35").Value = [double]"64.938446"
$ws.Range("O35").Value = [double]"0.05533777012737728"
$ws.Range("P35").Value = [double]"0.03782235733776705"
$ws.Range("Q35").Value = [double]"8.490690991425668"
$ws.Range("R35").Value = [double]"50.944145948554"
$ws.Range("S35").Value = [double]"3.009563789171603E-05"
$ws.Range("T35").Value = [double]"2.433927802545315E-05"
# Row 36
$ws.Range("G36").Value = [double]"0.2614996666666667"
$ws.Range("H36").Value = [double]"0.7844990000000001"
$ws.Range("I36").Value = [double]"0.0005438534625164957"
$ws.Range("J36").Value = [double]"0.0006435156277567465"
$ws.Range("M36").Value = [double]"52.09024033333333"
$ws.Range("N36").Value = [double]"156.270721"
$ws.Range("O36").Value = [double]"0.08877815602319267"
$ws.Range("P36").Value = [double]"0.09101737745760805"
$ws.Range("Q36").Value = [double]"13.62158048375322"
$ws.Range("R36").Value = [double]"122.594224353779"
$ws.Range("S36").Value = [double]"4.828230754904302E-05"
$ws.Range("T36").Value = [double]"5.857110479140539E-05"
# Row 37
$ws.Range("G37").Value = [double]"0.2614996666666667"
$ws.Range("H37").Value = [double]"0.7844990000000001"
$ws.Range("I37").Value = [double]"0.0005438534625164957"
$ws.Range("J37").Value = [double]"0.0006435156277567465"
$ws.Range("M37").Value = [double]"201.5345866666667"
$ws.Range("N37").Value = [double]"604.60376"
$ws.Range("O37").Value = [double]"0.3434783342267227"
$ws.Range("P37").Value = [double]"0.3521417722019025"
$ws.Range("Q37").Value = [double]"52.70122723513778"
$ws.Range("R37").Value = [double]"474.31104511624"
$ws.Range("S37").Value = [double]"0.0001868018813686013"
$ws.Range("T37").Value = [double]"0.0002266087335978805"

Write-Host "Updated all cells"
